$d = $word.ActiveDocument

# Remove the explicit "Justify" paragraph alignment (w:jc val="both") from
# every paragraph in the document, reverting them to the (unset/default)
# left alignment - this drops the <w:jc w:val="both"/> element entirely,
# and for paragraphs whose <w:pPr> only contained that element, the now
# empty <w:pPr> is dropped as well.
foreach ($p in $d.Paragraphs) {
    $p.Alignment = 0   # wdAlignParagraphLeft
}
